# CIVIC, NCI Thesaurus, MitelmanDB and Depmap update
$wb = $excel.ActiveWorkbook

$compounds = $wb.Worksheets.Item("compounds")
$biomarkers = $wb.Worksheets.Item("biomarkers")

# Mitelman Database row (row 3 of biomarkers): source_version v20240415 -> v20240715
$biomarkers.Range("E3").Value = "v20240715"

# PubChem row (row 5 of compounds): source_version v2023 -> v2024
$compounds.Range("E5").Value = "v2024"

# Add new DepMap row (row 5 of biomarkers)
$biomarkers.Range("F5").Value = "depmap"
$biomarkers.Range("F5").Style = "Normal"
$biomarkers.Range("C5").Value = "https://depmap.org"
$biomarkers.Range("D5").Value = "Tsherniak et al., Cell, 2017; 28753430"
$biomarkers.Range("E5").Value = "24Q2"
$biomarkers.Range("E5").Style = "Normal"
$biomarkers.Range("H5").Value = "https://depmap.org/portal/data_page/?tab=overview"
$biomarkers.Range("A5").Value = "DepMap"
$biomarkers.Range("B5").Value = "The Cancer Dependency Map"
$biomarkers.Range("G5").Value = "Free/open access"

# NCI Thesaurus row (row 3 of compounds): source_version 24.05d -> 24.07e
$compounds.Range("E3").Value = "24.07e"

# restore the selection/active cell state recorded in the saved view
$biomarkers.Range("B6").Select()
$compounds.Range("E3").Select()
